# Update "want to go" counts (column F) for a handful of events.
# These values are duplicated across the "展览" and "全部类型" sheets
# (same events), plus a single update on the "演出" sheet.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 (Exhibition) sheet
$wsExhibit.Range("F6").Value  = 1157
$wsExhibit.Range("F11").Value = 9501
$wsExhibit.Range("F22").Value = 2963
$wsExhibit.Range("F40").Value = 1470
$wsExhibit.Range("F47").Value = 734

# 演出 (Show) sheet
$wsShow.Range("F10").Value = 4

# 全部类型 (All types) sheet
$wsAll.Range("F6").Value  = 1157
$wsAll.Range("F8").Value  = 9501
$wsAll.Range("F19").Value = 2963
$wsAll.Range("F38").Value = 1470
$wsAll.Range("F47").Value = 734
$wsAll.Range("F49").Value = 4
